$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 18 (this shifts rows 18..169 down to 19..170,
# exactly matching the target diff: old row18 -> new row19, ..., old row169 -> new row170)
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with this week's record.
$ws.Range("A18").Value2 = 9
$ws.Range("B18").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C18").Value2 = "Metropolitana"
$ws.Range("D18").Value2 = 44490
$ws.Range("E18").Value2 = 13
$ws.Range("F18").Value2 = 300000001
$ws.Range("G18").Value2 = "Rabanito"
$ws.Range("H18").Value2 = "Sin especificar"
$ws.Range("I18").Value2 = "Primera"
$ws.Range("J18").Value2 = 7900
$ws.Range("K18").Value2 = 3000
$ws.Range("L18").Value2 = 4000
$ws.Range("M18").Value2 = 3494
$ws.Range("N18").Value2 = "`$/cien unidades (volumen en unidades)"
$ws.Range("O18").Value2 = "Provincia de Chacabuco"
$ws.Range("P18").Value2 = 35
$ws.Range("Q18").Value2 = 100
$ws.Range("R18").Value2 = "Hortaliza"
